$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (file1) separators first, to mirror shared-string insertion order
$ws.Range("E2").Value = "|"
$ws.Range("F2").Value = ";"

# New header cells
$ws.Range("E1").Value = "source_separator"
$ws.Range("F1").Value = "target_separator"
$ws.Range("G1").Value = "source_extension"
$ws.Range("H1").Value = "target_extension"

# Row 2 (file1) extensions
$ws.Range("G2").Value = "csv"
$ws.Range("H2").Value = "csv"

# Row 3 (file2)
$ws.Range("E3").Value = ";"
$ws.Range("F3").Value = ";"
$ws.Range("G3").Value = "csv"
$ws.Range("H3").Value = "csv"

# Column widths (closest achievable values given internal pixel-based rounding)
$ws.Range("E1").ColumnWidth = 15.833333333333334
$ws.Range("F1").ColumnWidth = 14.5
$ws.Range("G1").ColumnWidth = 14.666666666666666
$ws.Range("H1").ColumnWidth = 16.333333333333332

# Update selection to match target state
$ws.Range("C8").Select()
